$wb = $excel.ActiveWorkbook

# --- Overview sheet: status text updated from "Ready for handoff" to
#     "Handed back: in sync with en-US" for both locale columns (zh-cn / de-de)
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "Handed back: in sync with en-US"
$wsOverview.Range("F2").Value = "Handed back: in sync with en-US"
$wsOverview.Range("E3").Value = "Handed back: in sync with en-US"
$wsOverview.Range("F3").Value = "Handed back: in sync with en-US"

# --- zh-cn sheet: fill in "Latest Target File" (I) / "Latest Handback File" (J)
#     / "Latest Handback DateTime" (K) now that handback has happened, and
#     update the per-file Status column (C) to match the Overview sheet.
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("C2").Value = "Handed back: in sync with en-US"
$wsZh.Range("C3").Value = "Handed back: in sync with en-US"

$wsZh.Range("I2").Value = "a.md"
$wsZh.Hyperlinks.Add($wsZh.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/9378f215a8a0e13395494d3cf29b1eed34158dce/e2e/a.md", "", "", "a.md")
$wsZh.Range("J2").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf"
$wsZh.Range("K2").Value = "2016-09-02 20:43:41"

$wsZh.Range("I3").Value = "a.md"
$wsZh.Hyperlinks.Add($wsZh.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/9378f215a8a0e13395494d3cf29b1eed34158dce/e2e/a.md", "", "", "a.md")
$wsZh.Range("J3").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf"
$wsZh.Range("K3").Value = "2016-09-02 20:43:41"

# --- de-de sheet: same handback fields, but with de-de specific handback
#     filename/time.
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("C2").Value = "Handed back: in sync with en-US"
$wsDe.Range("C3").Value = "Handed back: in sync with en-US"

$wsDe.Range("I2").Value = "a.md"
$wsDe.Hyperlinks.Add($wsDe.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/9378f215a8a0e13395494d3cf29b1eed34158dce/e2e/a.md", "", "", "a.md")
$wsDe.Range("J2").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"
$wsDe.Range("K2").Value = "2016-09-02 20:43:48"

$wsDe.Range("I3").Value = "a.md"
$wsDe.Hyperlinks.Add($wsDe.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/9378f215a8a0e13395494d3cf29b1eed34158dce/e2e/a.md", "", "", "a.md")
$wsDe.Range("J3").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"
$wsDe.Range("K3").Value = "2016-09-02 20:43:48"
